# Auto-generated update of computed market-price / profit columns (H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @("H86", 2533.3333),
    @("I86", 1125),
    @("J86", 2727.5862),
    @("K86", 1125),
    @("L86", 2727.5862),
    @("M86", -2),
    @("N86", -4973.5862),
    @("H89", 2533.3333),
    @("I89", 1125),
    @("J89", 2727.5862),
    @("K89", 5625),
    @("L89", 13637.931),
    @("M89", -9),
    @("N89", -24869.931),
    @("H92", 1111.9286),
    @("I92", 873.8182),
    @("J92", 1985),
    @("K92", 873.8182),
    @("L92", 1985),
    @("M92", 374.1818),
    @("N92", -4481),
    @("H125", 13200.75),
    @("I125", 20520),
    @("J125", 1002),
    @("K125", 184680),
    @("L125", 9018),
    @("M125", -182220),
    @("N125", -13938),
    @("H137", 2063.5881),
    @("I137", 2152.0667),
    @("J137", 1400),
    @("K137", 6456.2001),
    @("L137", 4200),
    @("M137", -3906.2001),
    @("N137", -9300)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @("H45", 1778.4445),
    @("I45", 1818.3529),
    @("J45", 1100),
    @("K45", 1818.3529),
    @("L45", 1100),
    @("M45", -1441.3529),
    @("N45", -1854),
    @("H74", 6695.273),
    @("I74", 8988.615),
    @("K74", 8988.615),
    @("M74", -8114.615),
    @("H77", 6695.273),
    @("I77", 8988.615),
    @("K77", 44943.075),
    @("M77", -40575.075),
    @("H88", 3966.6667),
    @("I88", 3900),
    @("K88", 3900),
    @("M88", -3494),
    @("H91", 3966.6667),
    @("I91", 3900),
    @("K91", 3900),
    @("M91", -2496),
    @("H110", 735.0833),
    @("I110", 735.0833),
    @("J110", 0),
    @("K110", 735.0833),
    @("L110", 0),
    @("M110", 1309.9167),
    @("N110", $null),
    @("H122", 2188.7273),
    @("I122", 1799.6666),
    @("J122", 2655.6),
    @("K122", 5398.9998),
    @("L122", 7966.799999999999),
    @("M122", -2948.9998),
    @("N122", -12866.8)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @("H86", 2673.6924),
    @("I86", 2660.543),
    @("J86", 2788.75),
    @("K86", 2660.543),
    @("L86", 2788.75),
    @("M86", -1537.543),
    @("N86", -5034.75),
    @("H89", 2673.6924),
    @("I89", 2660.543),
    @("J89", 2788.75),
    @("K89", 13302.715),
    @("L89", 13943.75),
    @("M89", -7686.715),
    @("N89", -25175.75),
    @("H107", 1420),
    @("I107", 1413.3334),
    @("J107", 1433.3334),
    @("K107", 1413.3334),
    @("L107", 1433.3334),
    @("M107", 506.6666),
    @("N107", -5273.3334)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @("H31", 6469.68),
    @("I31", 8992.308000000001),
    @("K31", 8992.308000000001),
    @("M31", -8697.308000000001),
    @("H34", 6469.68),
    @("I34", 8992.308000000001),
    @("K34", 8992.308000000001),
    @("M34", -8790.308000000001),
    @("H96", 13600),
    @("J96", 13600),
    @("L96", 13600),
    @("N96", -19092),
    @("H105", 2658.889),
    @("H106", 196666.33),
    @("J106", 196666.33),
    @("L106", 196666.33),
    @("N106", -199190.33)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @("H55", 2816.3635),
    @("J55", 2998),
    @("L55", 8994),
    @("N55", -9348),
    @("H75", 1939.2858),
    @("I75", 1000),
    @("J75", 2095.8333),
    @("K75", 3000),
    @("L75", 6287.499899999999),
    @("M75", -2002),
    @("N75", -8283.499899999999),
    @("H78", 1939.2858),
    @("I78", 1000),
    @("J78", 2095.8333),
    @("K78", 9000),
    @("L78", 18862.4997),
    @("M78", -4008),
    @("N78", -28846.4997)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @("H122", 950),
    @("I122", 950),
    @("K122", 2850),
    @("M122", -400),
    @("H126", 0),
    @("I126", 0),
    @("J126", 0),
    @("K126", 0),
    @("L126", 0),
    @("M126", $null),
    @("N126", $null),
    @("H132", 28866.75),
    @("I132", 28234),
    @("K132", 84702),
    @("M132", -82172)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @("H7", 2710),
    @("I7", 2377.1428),
    @("J7", 3001.25),
    @("K7", 2377.1428),
    @("L7", 3001.25),
    @("M7", -2265.1428),
    @("N7", -3225.25),
    @("H126", 2710),
    @("I126", 2377.1428),
    @("J126", 3001.25),
    @("K126", 7131.428400000001),
    @("L126", 9003.75),
    @("M126", -4661.428400000001),
    @("N126", -13943.75)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @("H122", 359904.34),
    @("I122", 502151.16),
    @("J122", 4287.375),
    @("K122", 1506453.48),
    @("L122", 12862.125),
    @("M122", -1504003.48),
    @("N122", -17762.125)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
